$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price/volume data (GitHub Actions scheduled refresh).
# Price cells in column D are forced to text with a leading apostrophe when
# their new value would otherwise be auto-parsed as a number by Excel
# (which would strip meaningful trailing zeros / introduce float noise).
$ws.Range("D2").Value = "62.576.83"
$ws.Range("E2").Value = "  +4.38%  "
$ws.Range("D3").Value = "2.431.88"
$ws.Range("E3").Value = "  +5.40%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'558.29"
$ws.Range("E5").Value = "  +3.11%  "
$ws.Range("D6").Value = "'139.41"
$ws.Range("E6").Value = "  +7.53%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("E8").Value = "  +1.90%  "
$ws.Range("D9").Value = "2.430.92"
$ws.Range("E9").Value = "  +5.40%  "
$ws.Range("E10").Value = "  +3.84%  "
$ws.Range("E11").Value = "  +4.13%  "
$ws.Range("E12").Value = "  +0.49%  "
$ws.Range("E13").Value = "  +5.10%  "
$ws.Range("D14").Value = "'26.25"
$ws.Range("E14").Value = "  +12.48%  "
$ws.Range("D15").Value = "2.864.82"
$ws.Range("E15").Value = "  +5.39%  "
$ws.Range("D16").Value = "62.468.01"
$ws.Range("E16").Value = "  +4.25%  "
$ws.Range("E17").Value = "  +7.49%  "
$ws.Range("D18").Value = "2.432.87"
$ws.Range("E18").Value = "  +5.04%  "
$ws.Range("D19").Value = "'11.26"
$ws.Range("E19").Value = "  +7.38%  "
$ws.Range("D20").Value = "'347.12"
$ws.Range("E20").Value = "  +11.19%  "
$ws.Range("D21").Value = "'4.21"
$ws.Range("E21").Value = "  +3.21%  "
$ws.Range("D22").Value = "'6.83"
$ws.Range("E22").Value = "  +4.15%  "
$ws.Range("D23").Value = "'1.00"
$ws.Range("E23").Value = "  +0.03%  "
$ws.Range("D24").Value = "'5.55"
$ws.Range("E24").Value = "  -2.16%  "
$ws.Range("D25").Value = "'65.58"
$ws.Range("E25").Value = "  +2.99%  "
$ws.Range("E26").Value = "  +1.88%  "
$ws.Range("D27").Value = "'1.56"
$ws.Range("E27").Value = "  +16.05%  "
$ws.Range("D28").Value = "'0.998"
$ws.Range("E28").Value = "  -0.17%  "
$ws.Range("E29").Value = "  +5.44%  "
$ws.Range("E30").Value = "  +15.76%  "
$ws.Range("E31").Value = "  +5.80%  "
$ws.Range("E32").Value = "  +8.40%  "
$ws.Range("D33").Value = "'6.48"
$ws.Range("E33").Value = "  +11.30%  "
$ws.Range("D34").Value = "'172.63"
$ws.Range("E34").Value = "  +0.61%  "
$ws.Range("D35").Value = "'1.45"
$ws.Range("E35").Value = "  +7.50%  "
$ws.Range("E36").Value = "  +4.88%  "
$ws.Range("D37").Value = "'379.30"
$ws.Range("E37").Value = "  +20.09%  "
$ws.Range("D38").Value = "'18.61"
$ws.Range("E38").Value = "  +5.20%  "
$ws.Range("E39").Value = "  +11.83%  "
$ws.Range("E40").Value = "  +0.00%  "
$ws.Range("D41").Value = "'1.00"
$ws.Range("E41").Value = "  -0.18%  "
$ws.Range("E42").Value = "  +12.84%  "
$ws.Range("D43").Value = "'39.59"
$ws.Range("E43").Value = "  +3.90%  "
$ws.Range("D44").Value = "'144.79"
$ws.Range("E44").Value = "  +6.34%  "
$ws.Range("E45").Value = "  +7.47%  "
$ws.Range("D46").Value = "'20.83"
$ws.Range("E46").Value = "  +10.84%  "
$ws.Range("D47").Value = "'0.595"
$ws.Range("E47").Value = "  +4.42%  "
$ws.Range("B48").Value = "Hedera"
$ws.Range("C48").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D48").Value = "'0.0521"
$ws.Range("E48").Value = "  +6.44%  "
$ws.Range("B49").Value = "Stellar"
$ws.Range("C49").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D49").Value = "'0.0954"
$ws.Range("E49").Value = "  +1.71%  "
$ws.Range("E50").Value = "  +4.89%  "
$ws.Range("E51").Value = "  +6.61%  "
